$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two "category header" rows (5 and 8), shifting subsequent rows up.
# Delete the lower one first so row numbers for the other deletion stay valid.
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(5).Delete()

# Update the column-group sub-header in B2 (was "unnamed: 1_level_1") to "total"
$ws.Range("B2").Value = "total"
